# Update min_price (column D, numeric) and the matching price-string
# column (column F, text) for the rows whose prices changed.
#
# Each data row in this sheet stores the same price twice:
#   D<row> -> numeric value
#   F<row> -> the same value as text
# This script rewrites both columns for every affected row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    8   = 1449900
    9   = 767000
    33  = 1455000
    34  = 584900
    37  = 905000
    38  = 1545000
    40  = 1494990
    45  = 1859990
    50  = 1304400
    51  = 919000
    52  = 1007000
    53  = 615000
    54  = 1429000
    55  = 709000
    60  = 2252400
    61  = 1796940
    64  = 473000
    78  = 304400
    79  = 318500
    80  = 315500
    81  = 298500
    83  = 290900
    84  = 567900
    85  = 480000
    86  = 487000
    87  = 1052900
    88  = 1051900
    89  = 859900
    90  = 789900
    103 = 614000
    111 = 615000
    119 = 662400
    120 = 1534900
}

foreach ($row in $changes.Keys) {
    $newValue = $changes[$row]

    # D holds the true numeric price.
    $ws.Range("D$row").Value = $newValue

    # F holds the very same price again, but stored as text (it mirrors D
    # as a string in the source data). A leading apostrophe tells Excel to
    # keep the numeric-looking entry as text instead of re-parsing it back
    # into a number.
    $ws.Range("F$row").Value = "'" + [string]$newValue
}

$wb.Save()
